# This workbook's rows 4, 5 and 6 hold observation records that need to be
# cyclically shifted "up": row 4 receives row 5's data, row 5 receives row
# 6's data, and row 6 receives the original row 4's data (for the columns
# that actually carry per-observation data).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns whose values differ between rows 4/5/6 and therefore need to move.
$cols = @("A", "B", "E", "F", "G", "H", "P", "Q", "R", "S", "AI")

# Capture the original values for rows 4, 5 and 6 before overwriting anything.
# Value2 is used (rather than Value) so the raw underlying data is read/written.
$row4 = @{}
$row5 = @{}
$row6 = @{}
foreach ($col in $cols) {
    $row4[$col] = $ws.Range($col + "4").Value2
    $row5[$col] = $ws.Range($col + "5").Value2
    $row6[$col] = $ws.Range($col + "6").Value2
}

# Apply the cyclic rotation: row4 <- row5, row5 <- row6, row6 <- row4(original).
foreach ($col in $cols) {
    $ws.Range($col + "4").Value2 = $row5[$col]
    $ws.Range($col + "5").Value2 = $row6[$col]
    $ws.Range($col + "6").Value2 = $row4[$col]
}
